# Apply updated odds values to Sheet1, as captured by the source XML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    "O2"  = 1.36
    "P2"  = 3.2

    "M5"  = 1.05
    "O5"  = 1.29
    "Q5"  = 1.93
    "R5"  = 1.93

    "G6"  = 1.44
    "H6"  = 4.2
    "I6"  = 5.7
    "J6"  = 1.91
    "K6"  = 2.35
    "L6"  = 5.5
    "M6"  = 1.03
    "N6"  = 12.4
    "O6"  = 1.16
    "P6"  = 3.94
    "Q6"  = 1.62
    "R6"  = 2.05
    "S6"  = 1.3
    "T6"  = 3.22
    "U6"  = 1.84
    "V6"  = 1.92
    "W6"  = 6.5
    "X6"  = 6.2
    "Y6"  = 7
    "Z6"  = 8.25
    "AA6" = 9.5
    "AB6" = 19
    "AC6" = 12.5
    "AD6" = 7.3
    "AG6" = 350
    "AH6" = 14
    "AI6" = 29
    "AJ6" = 15
    "AK6" = 80
    "AO6" = 6.6
    "AP6" = 15.5
    "AQ6" = 18.5
    "AR6" = 45
    "AS6" = 200
    "AT6" = 3.1
    "AV6" = 70
    "AW6" = 7.4
    "AX6" = 32
    "BB6" = 450

    "U8"  = 1.67

    "G9"  = 1.5
    "H9"  = 4.1
    "I9"  = 6.25
    "J9"  = 2
    "K9"  = 2.4
    "M9"  = 1.04
    "N9"  = 13
    "Q9"  = 1.65
    "U9"  = 1.73
    "V9"  = 2
    "X9"  = 8
    "Z9"  = 11
    "AE9" = 15
    "AJ9" = 19
    "AN9" = 3.6
    "AO9" = 7.5
    "AU9" = 8
    "AW9" = 7.5
    "AX9" = 29
    "BA9" = 101

    "AT10" = 2.63

    "G13" = 1.39
}

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}
